$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New todo entry for the "Self-Assignable Roles" row.
$ws.Range("B5").Value = "add possibility for owners to make a default join role?"

# New "Track given interactions aswell" header column (E) on the todo table,
# styled like the other header cells (A3:D3 use the "Gut"/Good style).
$ws.Range("E3").Value = "Track given interactions aswell"
$ws.Range("E3").Style = "Gut"

# Give column E a custom width matching the new table column.
$ws.Columns("E").ColumnWidth = 45.6

# Highlight the "Info" module row: label cell gets the new "Input" style,
# its two todo cells reuse the existing "Gut" (Good) style.
$ws.Range("A19").Style = "Input"
$ws.Range("B19").Style = "Gut"
$ws.Range("C19").Style = "Gut"

# Move the active selection to C9 (matches the saved view state).
$ws.Range("C9").Select()
